# Add 2022-Q3 data
# 1) Insert a new "2022-Q3" worksheet right after "总计", populated as a copy of
#    the "2022-Q2" sheet (so header/column formatting matches), then overwrite
#    its data with the new Q3 numbers.
# 2) Update the "总计" (summary) sheet to add a row for 2022-Q3 and shift the
#    existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q3" sheet
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# Copy "2022-Q2" right after "总计" -- gives us the same header row / styling,
# which we then overwrite with the new quarter's values.
$wsQ2.Copy($null, $wsTotal)
$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "2022-Q3"

# Row 2 -- 013466 / 博时智选量化多因子股票C
$wsNew.Range("B2").NumberFormat = "@"
$wsNew.Range("B2").Value = "013466"
$wsNew.Range("C2").Value = "博时智选量化多因子股票C"
$wsNew.Range("D2").NumberFormat = "@"
$wsNew.Range("D2").Value = "2.28"
$wsNew.Range("E2").NumberFormat = "@"
$wsNew.Range("E2").Value = "92.38"
$wsNew.Range("F2").NumberFormat = "@"
$wsNew.Range("F2").Value = "1.06"
$wsNew.Range("G2").NumberFormat = "@"
$wsNew.Range("G2").Value = "0.0242"
$wsNew.Range("H2").Value = 8

# Row 3 -- 013465 / 博时智选量化多因子股票A
$wsNew.Range("B3").NumberFormat = "@"
$wsNew.Range("B3").Value = "013465"
$wsNew.Range("C3").Value = "博时智选量化多因子股票A"
$wsNew.Range("D3").NumberFormat = "@"
$wsNew.Range("D3").Value = "0.49"
$wsNew.Range("E3").NumberFormat = "@"
$wsNew.Range("E3").Value = "92.38"
$wsNew.Range("F3").NumberFormat = "@"
$wsNew.Range("F3").Value = "1.06"
$wsNew.Range("G3").NumberFormat = "@"
$wsNew.Range("G3").Value = "0.0052"
$wsNew.Range("H3").Value = 8

# ---------------------------------------------------------------------------
# Step 2: update the "总计" sheet -- insert a row for 2022-Q3 and shift the
# remaining quarters (and their running index in column A) down by one.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")
$ws1.Rows.Item(2).Insert()

# Fix up formatting on the newly inserted row: column A should carry the same
# style as the other index cells, columns B:D should have no special style
# (Insert() copies the header row's style into the new blank row by default).
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

$ws1.Range("B3:D3").Copy()
$ws1.Range("B2:D2").PasteSpecial(-4122)

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 0.03

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 2
$ws1.Range("D3").Value = 0.04

$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "2022-Q1"
$ws1.Range("C4").Value = 7
$ws1.Range("D4").Value = 0.79

$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "2021-Q4"
$ws1.Range("C5").Value = 2
$ws1.Range("D5").Value = 0.06
